$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sheet1")

$ws.Range("C3").Value = "selidik"
$ws.Range("C4").Value = "bunuh"
$ws.Range("C37").Value = "selid"
$ws.Range("C53").Value = "duduk"
$ws.Range("C104").Value = "bakar"
$ws.Range("C105").Value = "perkosa"
